$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as row 116, pushing the
# existing rows 116:190 down to 117:191 (dimension grows from R190 to R191).
$ws.Rows.Item(116).Insert()

$ws.Range("A116").Value = 8
$ws.Range("B116").Value = "Terminal La Palmera de La Serena"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value2 = 45090
$ws.Range("E116").Value = 4
$ws.Range("F116").Value = 100112052
$ws.Range("G116").Value = "Albahaca"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 700
$ws.Range("K116").Value = 2800
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 2900
$ws.Range("N116").Value = "`$/paquete"
$ws.Range("O116").Value = "Región de Arica y Parinacota"
$ws.Range("P116").Value = 2900
$ws.Range("Q116").Value = 1
$ws.Range("R116").Value = "Hortaliza"
